$d = $word.ActiveDocument

# Locate and remove the "....By Siddhesh Shivtarkar" byline paragraph
# (including its own paragraph mark) as requested by the author.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*By Siddhesh Shivtarkar*") {
        $rng = $p.Range
        $rng.Delete()
        # Word marks the location of the last edit with the "_GoBack" bookmark.
        $rng.Collapse(1)
        $d.Bookmarks.Add("_GoBack", $rng)
        break
    }
}
